$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 113, shifting rows 113:217 down to 114:218
$ws.Rows(113).Insert()

# Populate the newly inserted row 113 with the new record
$ws.Cells.Item(113, 1).Value = 9
$ws.Cells.Item(113, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(113, 3).Value = "Metropolitana"
$ws.Cells.Item(113, 4).Value = 44494
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = 100112032
$ws.Cells.Item(113, 7).Value = "Zapallo italiano"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 130
$ws.Cells.Item(113, 11).Value = 10000
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 13).Value = 10000
$ws.Cells.Item(113, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 167
$ws.Cells.Item(113, 17).Value = 60
$ws.Cells.Item(113, 18).Value = "Hortaliza"
